$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New prediction rows appended by the automatic tracker sync (rows 30-42).
# Helper that writes a plain text value to a cell without letting Excel
# auto-convert look-alike dates/percentages into numbers, and without
# leaving a residual number format style behind on the cell.
function Set-TextCell($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$newRows = @(
    @("30", '2025-08-21', 'Ligue 1', 'paris saint germain', 'angers', 'Away Win', '98.41%', 21, '1945.95%', 6.2, 0.05, 0.9833078686563578, 'Pending', '', '', '', '', ''),
    @("31", '2025-08-21', 'La Liga', 'real betis', 'alaves', 'Away Win', '48.73%', 4.6, '121.91%', 4.3, 0.03448745631990931, 0.344874563199093, 'Pending', '', '', '', '', ''),
    @("32", '2025-08-21', 'Segunda División', 'eibar', 'granada cf', 'Away Win', '44.20%', 3.6, '57.52%', 2.8, 0.02273606502137256, 0.2273606502137256, 'Pending', '', '', '', '', ''),
    @("33", '2025-08-21', 'Segunda División', 'leganes', 'cadiz', 'Away Win', '47.78%', 4.5, '112.85%', 4.1, 0.03285799271357884, 0.3285799271357884, 'Pending', '', '', '', '', ''),
    @("34", '2025-08-21', 'Premier League', 'west ham', 'chelsea', 'Home Win', '59.34%', 4.75, '179.07%', 6, 0.04850299398501735, 0.4850299398501735, 'Pending', '', '', '', '', ''),
    @("35", '2025-08-21', 'Eliteserien', 'kfum oslo', 'ham-kam', 'Away Win', '53.43%', 5.25, '177.69%', 5.3, 0.04247001729525227, 0.4247001729525227, 'Pending', '', '', '', '', ''),
    @("36", '2025-08-21', 'Ekstraklasa', 'radomiak radom', 'nieciecza', 'Away Win', '43.47%', 3.6, '54.94%', 2.7, 0.02173134634323324, 0.2173134634323323, 'Pending', '', '', '', '', ''),
    @("37", '2025-08-21', 'Superliga', 'vllaznia shkodër', 'egnatia rrogozhinë', 'Away Win', '40.15%', 3.25, '29.19%', 1.7, 0.01355307647359828, 0.1355307647359828, 'Pending', '', '', '', '', ''),
    @("38", '2025-08-21', '2. Bundesliga', 'sv elversberg', '1. fc kaiserslautern', 'Away Win', '40.93%', 3.25, '31.70%', 1.8, 0.01467980316853604, 0.1467980316853604, 'Pending', '', '', '', '', ''),
    @("39", '2025-08-21', 'League', 'sama al sarhan', 'al salt', 'Home Win', '53.30%', 4.2, '121.61%', 4.8, 0.03870251266414631, 0.387025126641463, 'Pending', '', '', '', '', ''),
    @("40", '2025-08-21', 'Liga Profesional Argentina', 'tigre', 'independ. rivadavia', 'Away Win', '42.39%', 3.8, '59.45%', 2.7, 0.02180870899900153, 0.2180870899900152, 'Pending', '', '', '', '', ''),
    @("41", '2025-08-21', 'Primera B', 'real cartagena', 'tigres fc', 'Away Win', '69.06%', 8, '446.92%', 6.2, 0.05, 0.6463526766868029, 'Pending', '', '', '', '', ''),
    @("42", '2025-08-21', '1. Division', 'aalborg', 'hvidovre', 'Away Win', '42.99%', 3.4, '44.69%', 2.4, 0.01922983319639007, 0.1922983319639007, 'Pending', '', '', '', '', '')
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")
$numCols  = @("G","I","J","K")

foreach ($row in $newRows) {
    $r = $row[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $val = $row[$i + 1]
        $cellRef = "$col$r"
        if ($numCols -contains $col) {
            $ws.Range($cellRef).Value = [double]$val
        } else {
            Set-TextCell $cellRef $val
        }
    }
}
